$d = $word.ActiveDocument

# --- Step 1: merge the two runs of paragraph 9 into a single run -----------
# "9. " + "infinite loop with key Ctrl + C." -> "9. infinite loop with key Ctrl + C."
$d.Content.Find.Execute("9. infinite loop with key Ctrl + C.", $false, $false,
                         $false, $false, $false, $true, 1, $false,
                         "9. infinite loop with key Ctrl + C.", 2)

# --- Step 2: insert a new paragraph "10. ..." right after paragraph 9 ------
$para9 = $d.Paragraphs.Item(10)
$para9.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(11)
$nr = $newPara.Range
$nr.MoveEnd(1, -1)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">10. </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
      <w:color w:val="040C28"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t>Break statement stops the entire process of the loop.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
      <w:color w:val="202124"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
      <w:color w:val="040C28"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t>Continue statement only stops the current iteration of the loop</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
      <w:color w:val="202124"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
'@

$nr.InsertXML($xml)
